# 10th - MB for single stock and added new group
#
# The weekly MarketBeat rank tracker gets 2 new date snapshots (Jun_26 and
# Jun_27) inserted as the newest columns (pushing the existing Jun_17 /
# Jun_15 / Jun_13 / Jun_10 columns to the right), each new data cell
# defaulting to "UN" (unchanged), and a new ratings-source group
# (Benchmark, Evercore ISI) appended as two new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns in front of the existing date columns (old B:D -> E:G,
# old E -> H). Column B will become the newest "Jun_27" snapshot; C and D
# are the "Jun_26" snapshot (duplicated across both columns, matching the
# source report).
$ws.Columns("B:D").Insert()

# New column headers (row 1).
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Default all the newly inserted data cells (rows 2-27, columns B:D) to
# "UN" (unchanged) same as every other snapshot column.
$ws.Range("B2:D27").Value = "UN"

# New group of ratings sources appended at the bottom.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28:D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29:D29").Value = "UN"

# Latest snapshot header, set last.
$ws.Range("B1").Value = "Jun_27"
